# Update latest output (run 160)
$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": recompute cost / unit-cost columns -----------------
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E2").Value = 384.5930985000001
$schedule.Range("F2").Value = 8.478683829365082

$schedule.Range("E3").Value = 433.27181625
$schedule.Range("F3").Value = 28.65554340277778

# --- Sheet "Detailed": refreshed price forecast / historical rows --------
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B10").Value = 78
$detailed.Range("C10").Value = "historical"

$detailed.Range("B11").Value = 79.61741000000001
$detailed.Range("C11").Value = "historical"

$detailed.Range("B12").Value = 78

$detailed.Range("B13").Value = 80.96411000000001

$detailed.Range("B14").Value = 76.5692

$detailed.Range("B15").Value = 59.35883

$detailed.Range("B16").Value = 36.06

$detailed.Range("B18").Value = -5.19167

$detailed.Range("B19").Value = -6.21833

$detailed.Range("B20").Value = -7.90661

$detailed.Range("B21").Value = -7.98427

$detailed.Range("B22").Value = -7.75328

$detailed.Range("B23").Value = -7.48607

$detailed.Range("B24").Value = -7.35259

$detailed.Range("B25").Value = -0.37997

$detailed.Range("B26").Value = -6.8

$detailed.Range("B28").Value = -7.86135

$detailed.Range("B29").Value = -5.99525

$detailed.Range("B30").Value = -4.82645

$detailed.Range("B31").Value = -0.8763300000000001

$detailed.Range("B37").Value = -7.34703

$detailed.Range("B38").Value = -1.15013

$detailed.Range("B39").Value = 3.06249

$detailed.Range("B40").Value = 30.51578

$detailed.Range("B41").Value = 51.4753

$detailed.Range("B42").Value = 53.90468

$detailed.Range("B43").Value = 58.02287

$detailed.Range("B45").Value = 57.03541

$detailed.Range("B46").Value = 47.42917
